$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the hard coded travel dates / swap rows so data isn't hard-coded
$ws.Range("C2").Value = "15-04-2024"

$ws.Range("A3").Value = "DEL"
$ws.Range("B3").Value = "CHN"
$ws.Range("C3").Value = "20-03-2024"

$ws.Range("C4").Value = "25-02-2024"

$ws.Range("A5").Value = "DEL"
$ws.Range("B5").Value = "CHN"
$ws.Range("C5").Value = "13-06-2024"

# Update selection to A3
$ws.Range("A3").Select()
